# LOQ4066.xlsx edit: add "Objectives" text, insert new "Programa resumido" text,
# full "Programa" syllabus text, and full "Bibliografia" text; shift everything
# below "Programa resumido:" down by one row to make room.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert one new (blank) row right below the existing "Programa resumido:"/
#    "Semestral" row (row 14), pushing every row from 14 downward down by one.
#    Excel auto-copies A/B/C formatting (styles 1/2/3) from row 13 into the new
#    blank row 14, and keeps row/column formatting for everything pushed down.
# ---------------------------------------------------------------------------
$ws.Rows.Item(14).Insert()

# ---------------------------------------------------------------------------
# 2) New "Objetivos:" (row 10) body text, in both the Portuguese (B) and the
#    "changed" (C) column.
# ---------------------------------------------------------------------------
$objetivos = 'Apresentar os princípios fundamentais envolvidos nas operações relacionadas a sistemas particulados, de forma a permitir a análise de desempenho dos equipamentos que lidam com estes sistemas.'
$ws.Range("B10").Value = $objetivos
$ws.Range("C10").Value = $objetivos

# ---------------------------------------------------------------------------
# 3) Old row 13 ("Programa resumido:" / "Semestral") becomes the
#    "Docentes responsáveis:" value row: clear the A-label, and replace the
#    B/C text with the professor's name (moved down from where it used to be).
# ---------------------------------------------------------------------------
$ws.Range("A13").ClearContents()
$docente = '787307 - Luis Fernando Figueiredo Faria'
$ws.Range("B13").Value = $docente
$ws.Range("C13").Value = $docente
$ws.Rows.Item(13).EntireRow.AutoFit()

# ---------------------------------------------------------------------------
# 4) The newly inserted row 14 becomes the new "Programa resumido:" row, with
#    a brand-new summary paragraph in B/C.
# ---------------------------------------------------------------------------
$ws.Range("A14").Value = "Programa resumido:"
$resumido = 'Fundamentos e caracterização de partículas e sistemas particulados. Dinâmica da interação sólido-fluido. Aplicações em sistemas diluídos: elutriação, câmara de poeira, ciclones, centrífugas e hidrociclones. Aplicações em sistemas concentrados: escoamento monofásico em meios porosos, filtração sólido-líquido, sedimentação, fluidização, transporte pneumático e hidráulico de partículas.'
$ws.Range("B14").Value = $resumido
$ws.Range("C14").Value = $resumido
$ws.Rows.Item(14).RowHeight = 60

# ---------------------------------------------------------------------------
# 5) "Programa:" (now row 16) gets the full syllabus text in B/C.
# ---------------------------------------------------------------------------
$programa = @'
1. Caracterização de partículas e sistemas particulados: noções de amostragem; diâmetros de esferas equivalentes e  diâmetros estatísticos; esfericidade; análise granulométrica, frequência simples e acumuladas; modelos de distribuição de tamanhos.   
2. Interação sólido-fluido: Dinâmica e análise dimensional do sistema partícula - fluido infinito:velocidade terminal; lei de Stokes; correlação entre coeficiente de arraste e número de Reynolds para esferas; efeito da forma das partículas; efeito de paredes; efeito de população; efeito de deslizamento. 
3. Aplicações em sistemas diluídos: separação sólido-sólido por elutriação; partículas equitombantes e razão de sedimentação; eficiências globais e individuais de coleta; diâmetro de corte; separação sólido-gás com câmaras de poeira e ciclones; separação sólido-líquido com centrífugas e hidrociclones.
4. Aplicações em  sistemas concentrados: escoamento monofásico em meios porosos; separação sólido-líquido por filtração em superfície; auxiliares de filtração; estudo detalhado dos filtros prensa e de tambor rotativo; separação sólido-líquido por sedimentação; leitos fluidizados a gás e a líquido; curva característica e histerese de fluidização; previsão das velocidades mínima e máxima de fluidização; transporte pneumático de partículas; velocidade de deslizamento; transporte hidráulico de partículas; velocidade de salto.
'@
$programa = $programa.TrimEnd("`r", "`n")
$ws.Range("B16").Value = $programa
$ws.Range("C16").Value = $programa

# ---------------------------------------------------------------------------
# 6) "Método:" (now row 19) gets the evaluation-method text in B/C (same text
#    that used to sit under "Critério:").
# ---------------------------------------------------------------------------
$metodo = 'Participação em sala de aula, preparação e apresentação de trabalhos e provas escritas.'
$ws.Range("B19").Value = $metodo
$ws.Range("C19").Value = $metodo

# ---------------------------------------------------------------------------
# 7) "Critério:" (now row 20) gets the "Média Final" grading-formula text
#    (moved up from where "Norma de recuperação:" used to carry it).
# ---------------------------------------------------------------------------
$criterio = @'
Média Final = (Prova1 + Prova2 + Nota de Trabalho) /3
Média final mínima de aprovação = 5,0
'@
$criterio = $criterio.TrimEnd("`r", "`n")
$ws.Range("B20").Value = $criterio
$ws.Range("C20").Value = $criterio

# ---------------------------------------------------------------------------
# 8) "Norma de recuperação:" (now row 21) gets the retake-exam formula text
#    (moved up from where "Bibliografia:" used to carry it).
# ---------------------------------------------------------------------------
$norma = '(Prova escrita + Média Final)/2         Nota Final mínima para aprovação= 5,0'
$ws.Range("B21").Value = $norma
$ws.Range("C21").Value = $norma

# ---------------------------------------------------------------------------
# 9) "Bibliografia:" (now row 22) gets the full reading list.
# ---------------------------------------------------------------------------
$bibliografia = @'
1. PERRY, R.H.; GREEN, D.W.; MALONEY, J.O. (Eds.). Perrys Chemical Engineers Handbook. New York : McGraw-Hill, 1997.
2. MASSARANI, G. Fluidodinâmica em Sistemas Particulados. 2. ed. RJ: E-Papers, 2002.
3. SVAROVSKY, L. Solid-Liquid Separation. 3. ed. LondonBoston: Butterworths, 1990.
4. RUSHTON, A.; WARD, A.S.; HOLDICH, R.G. Solid-Liquid Filtration and Separation Technology. Weinheim:  VCH, 1996.
5. COULSON, J.M.; RICHARDSON, J.F. Chemical Engineering. 5th. ed. Londres: Pergamon Press,1996. Vol. 2.
6. ALLEN, T. Particle Size Measurement. 5th. ed. Londres: Chapman & Hall, 1997. Vol 1 e 2.
'@
$bibliografia = $bibliografia.TrimEnd("`r", "`n")
$ws.Range("B22").Value = $bibliografia
$ws.Range("C22").Value = $bibliografia
